$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "97.34") must be
# kept as TEXT (matching the original inlineStr cell type), so force a
# text number format before assigning the value.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D16", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D30", "D31", "D33", "D34", "D39", "D40", "D42", "D43", "D44", "D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.963.94"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.238.59"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "97.34"
$ws.Range("E5").Value = "  +17.22%  "
$ws.Range("D6").Value = "272.49"
$ws.Range("E6").Value = "  +5.34%  "
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").Value = "  +5.74%  "
$ws.Range("D10").Value = "47.89"
$ws.Range("E10").Value = "  +7.99%  "
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "8.29"
$ws.Range("E12").Value = "  +15.64%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "15.28"
$ws.Range("E14").Value = "  +6.38%  "
$ws.Range("D15").Value = "2.574.12"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "0.826"
$ws.Range("E16").Value = "  +5.56%  "
$ws.Range("D17").Value = "2.244.74"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "43.940.92"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +4.69%  "
$ws.Range("D21").Value = "70.89"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").Value = "234.14"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "9.38"
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "11.41"
$ws.Range("E26").Value = "  +7.10%  "
$ws.Range("E27").Value = "  +12.05%  "
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "39.49"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "172.91"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  +6.08%  "
$ws.Range("D33").Value = "21.07"
$ws.Range("E33").Value = "  +3.33%  "
$ws.Range("D34").Value = "5.59"
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "3.59"
$ws.Range("D40").Value = "0.249"
$ws.Range("E40").Value = "  +25.00%  "
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").Value = "12.46"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "61.95"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").Value = "100.40"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "2.456.79"
$ws.Range("E51").Value = "  +1.60%  "
